$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to match the player's name
$ws.Name = "Lalit Yadav"

# Header row (now 13 columns, with a new "matchNo" column inserted before teamName)
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Full set of innings data rows:
# matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result
$rows = @(
    @("41st","Delhi Capitals","Lalit Yadav","lbw b Narine","0","3","0","0","0.00","Kolkata Knight Riders","Sharjah","September 28","KKR won by 3 wickets (with 10 balls remaining)"),
    @("11th","Delhi Capitals","Lalit Yadav","","12","6","2","0","200.00","Punjab Kings","Wankhede","April 18","Capitals won by 6 wickets (with 10 balls remaining)"),
    @("36th","Delhi Capitals","Lalit Yadav","","14","15","1","0","93.33","Rajasthan Royals","Abu Dhabi","September 25","Capitals won by 33 runs"),
    @("13th","Delhi Capitals","Lalit Yadav","","22","25","1","0","88.00","Mumbai Indians","Chennai","April 20","Capitals won by 6 wickets (with 5 balls remaining)"),
    @("7th","Delhi Capitals","Lalit Yadav","c Tewatia b Morris","20","24","3","0","83.33","Rajasthan Royals","Wankhede","April 15","Royals won by 3 wickets (with 2 balls remaining)")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $val = $rowData[$c]

        if ($val -eq "") {
            # No dismissal info recorded for this innings (source keeps an
            # explicit, empty text cell rather than leaving it blank).
            $cell.Formula = "=""""" 
        } elseif ($val -match '^-?\d+(\.\d+)?$') {
            # Every value in this sheet is stored as text (even the
            # numeric-looking ones, e.g. "0.00" strike-rates or ball
            # counts) - format as Text first so Excel keeps the exact
            # string instead of coercing it to a number.
            $cell.NumberFormat = "@"
            $cell.Value = $val
        } else {
            $cell.Value = $val
        }
    }
}
